$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-115). Update it to the new "changed" date serial 45188
# (2023-09-19), matching the source diff.
$ws.Range("C2:C115").Value = 45188
